# Apply the "Fixed update to excel issue" change:
# 1. Rename the "Requested quantity" headers on the existing sheets.
# 2. Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- 1. Rename headers -------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" worksheet after the last sheet -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(44990.99999999999, 5, 5.000000109113884, 5.000000109406614),
    @(45137.99999999999, 1, 1.000000109323564, 1.000000109609726),
    @(45144.99999999999, 1, 0.8095239187258656, 0.8095239190212431),
    @(45151.99999999999, 1, 0.6190477281248392, 0.6190477284466256),
    @(45158.99999999999, 0, 0.4285715377840956, 0.4285715381451593),
    @(45165.99999999999, 0, 0.2380953471875631, 0.2380953475797541),
    @(45172.99999999999, 0, 0.04761915678249234, 0.04761915747260548),
    @(45179.99999999999, 0, -0.1428570344108191, -0.1428570321288324),
    @(45186.99999999999, 0, -0.3333332258860666, -0.3333332217845986),
    @(45193.99999999999, 0, -0.5238094171943806, -0.5238094107601413)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}

# --- 3. Match formatting of the existing sheets -------------------------
# Header style (bold, centered, bordered) - copy from the existing header cell
# so the same cell style gets reused instead of creating a brand new one.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

# Date column style (yyyy-mm-dd hh:mm:ss number format) - copy from the
# existing date cell for the same reason.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)  # xlPasteFormats

